# Auto-generated Excel COM-interop script to apply market-data/profit updates
# across the Titan_Profits workbook sheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2136.8
$ws.Range("J40").Value = 2000.4
$ws.Range("L40").Value = 2000.4
$ws.Range("N40").Value = -2350.4

$ws.Range("H64").Value = 348557.28
$ws.Range("I64").Value = 502983.1
$ws.Range("J64").Value = 5388.778
$ws.Range("K64").Value = 502983.1
$ws.Range("L64").Value = 5388.778
$ws.Range("M64").Value = -502735.1
$ws.Range("N64").Value = -5884.778

$ws.Range("H67").Value = 348557.28
$ws.Range("I67").Value = 502983.1
$ws.Range("J67").Value = 5388.778
$ws.Range("K67").Value = 502983.1
$ws.Range("L67").Value = 5388.778
$ws.Range("M67").Value = -502125.1
$ws.Range("N67").Value = -7104.778

$ws.Range("H133").Value = 11857
$ws.Range("J133").Value = 11857
$ws.Range("L133").Value = 11857
$ws.Range("N133").Value = -21977

$ws.Range("H134").Value = 70653.336
$ws.Range("J134").Value = 70653.336
$ws.Range("L134").Value = 70653.336
$ws.Range("N134").Value = -80793.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18293.146
$ws.Range("I32").Value = 3406.4736
$ws.Range("J32").Value = 95433.17999999999
$ws.Range("K32").Value = 3406.4736
$ws.Range("L32").Value = 95433.17999999999
$ws.Range("M32").Value = -3119.4736
$ws.Range("N32").Value = -96007.17999999999

$ws.Range("H74").Value = 5682.2583
$ws.Range("I74").Value = 1685.3334
$ws.Range("K74").Value = 1685.3334
$ws.Range("M74").Value = -811.3334

$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25676

$ws.Range("H77").Value = 5682.2583
$ws.Range("I77").Value = 1685.3334
$ws.Range("K77").Value = 8426.666999999999
$ws.Range("M77").Value = -4058.666999999999

$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27340

$ws.Range("H132").Value = 3241.975
$ws.Range("I132").Value = 2828.818
$ws.Range("J132").Value = 5189.7144
$ws.Range("K132").Value = 8486.454000000002
$ws.Range("L132").Value = 15569.1432
$ws.Range("M132").Value = -5956.454000000002
$ws.Range("N132").Value = -20629.1432

$ws.Range("H133").Value = 52199.8
$ws.Range("J133").Value = 52199.8
$ws.Range("L133").Value = 52199.8
$ws.Range("N133").Value = -57259.8

$ws.Range("H139").Value = 40734.57
$ws.Range("J139").Value = 40734.57
$ws.Range("L139").Value = 40734.57
$ws.Range("N139").Value = -51014.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 49950
$ws.Range("J59").Value = 49950
$ws.Range("L59").Value = 49950
$ws.Range("N59").Value = -51644

$ws.Range("H100").Value = 22475
$ws.Range("J100").Value = 22475
$ws.Range("L100").Value = 22475
$ws.Range("N100").Value = -24639

$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -55120

$ws.Range("H134").Value = 21278894
$ws.Range("I134").Value = 25642816
$ws.Range("J134").Value = 4776.5
$ws.Range("K134").Value = 76928448
$ws.Range("L134").Value = 14329.5
$ws.Range("M134").Value = -76925913
$ws.Range("N134").Value = -19399.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 500000
$ws.Range("I16").Value = 500000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 500000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -499713
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 1518.4777
$ws.Range("I31").Value = 906.93616
$ws.Range("J31").Value = 2955.6
$ws.Range("K31").Value = 906.93616
$ws.Range("L31").Value = 2955.6
$ws.Range("M31").Value = -611.93616
$ws.Range("N31").Value = -3545.6

$ws.Range("H34").Value = 1518.4777
$ws.Range("I34").Value = 906.93616
$ws.Range("J34").Value = 2955.6
$ws.Range("K34").Value = 906.93616
$ws.Range("L34").Value = 2955.6
$ws.Range("M34").Value = -704.93616
$ws.Range("N34").Value = -3359.6

$ws.Range("H113").Value = 500000
$ws.Range("I113").Value = 500000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 500000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -497830
$ws.Range("N113").ClearContents()

$ws.Range("H119").Value = 44630.5
$ws.Range("J119").Value = 44630.5
$ws.Range("L119").Value = 44630.5
$ws.Range("N119").Value = -54306.5

$ws.Range("H134").Value = 2435.1714
$ws.Range("I134").Value = 1651.2413
$ws.Range("J134").Value = 6224.1665
$ws.Range("K134").Value = 4953.7239
$ws.Range("L134").Value = 18672.4995
$ws.Range("M134").Value = -2418.7239
$ws.Range("N134").Value = -23742.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H107").Value = 486.6207
$ws.Range("I107").Value = 462.375
$ws.Range("J107").Value = 516.46155
$ws.Range("K107").Value = 1387.125
$ws.Range("L107").Value = 1549.38465
$ws.Range("M107").Value = 532.875
$ws.Range("N107").Value = -5389.38465

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1503.5
$ws.Range("I9").Value = 1004.6667
$ws.Range("K9").Value = 1004.6667
$ws.Range("M9").Value = -834.6667

$ws.Range("H132").Value = 2975.6394
$ws.Range("I132").Value = 2934
$ws.Range("J132").Value = 3049.4546
$ws.Range("K132").Value = 8802
$ws.Range("L132").Value = 9148.363799999999
$ws.Range("M132").Value = -6272
$ws.Range("N132").Value = -14208.3638

$ws.Range("H135").Value = 500022900
$ws.Range("J135").Value = 500022900
$ws.Range("L135").Value = 500022900
$ws.Range("N135").Value = -500033040

$ws.Range("H138").Value = 72080
$ws.Range("J138").Value = 72080
$ws.Range("L138").Value = 72080
$ws.Range("N138").Value = -82360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1398
$ws.Range("J46").Value = 1663.3334
$ws.Range("L46").Value = 1663.3334
$ws.Range("N46").Value = -2039.3334

$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622

$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1133.3334
$ws.Range("I96").Value = 1400
$ws.Range("K96").Value = 1400
$ws.Range("M96").Value = -27

$ws.Range("H132").Value = 2182.0344
$ws.Range("I132").Value = 2194.9119
$ws.Range("J132").Value = 2135.9473
$ws.Range("K132").Value = 6584.7357
$ws.Range("L132").Value = 6407.841899999999
$ws.Range("M132").Value = -4054.7357
$ws.Range("N132").Value = -11467.8419

$ws.Range("H136").Value = 43837.957
$ws.Range("I136").Value = 101471.1
$ws.Range("J136").Value = 2671.4285
$ws.Range("K136").Value = 304413.3
$ws.Range("L136").Value = 8014.2855
$ws.Range("M136").Value = -301863.3
$ws.Range("N136").Value = -13114.2855
